$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-03-01 05:18:33"
$ws.Range("N2").Value = "-1.9 °C 4:56 TU"
$ws.Range("O2").Value = "-0.8 °C"
$ws.Range("E3").Value = "2026-03-01 05:18:35"
$ws.Range("L3").Value = "17.6 km/h - 120º 4:47 TU"
$ws.Range("N3").Value = "-4.2 °C 4:59 TU"
$ws.Range("O3").Value = "-3.6 °C"
$ws.Range("E4").Value = "2026-03-01 05:18:37"
$ws.Range("E5").Value = "2026-03-01 05:18:40"
$ws.Range("N5").Value = "-4.5 °C 4:57 TU"
$ws.Range("O5").Value = "-3.7 °C"
$ws.Range("E6").Value = "2026-03-01 05:18:43"
$ws.Range("H6").Value = "88%"
$ws.Range("N6").Value = "8.4 °C 4:32 TU"
$ws.Range("O6").Value = "9.3 °C"
$ws.Range("E7").Value = "2026-03-01 05:18:46"
$ws.Range("H7").Value = "74%"
$ws.Range("N7").Value = "12.9 °C 4:48 TU"
$ws.Range("E8").Value = "2026-03-01 05:18:48"
$ws.Range("N8").Value = "9.2 °C 4:30 TU"
$ws.Range("E9").Value = "2026-03-01 05:18:51"
$ws.Range("E10").Value = "2026-03-01 05:18:53"
$ws.Range("N10").Value = "4.7 °C 4:50 TU"
$ws.Range("O10").Value = "6.5 °C"
$ws.Range("E11").Value = "2026-03-01 05:18:56"
$ws.Range("E12").Value = "2026-03-01 05:18:58"
$ws.Range("H12").Value = "72%"
$ws.Range("N12").Value = "9.0 °C 4:54 TU"
$ws.Range("O12").Value = "10.5 °C"
$ws.Range("E13").Value = "2026-03-01 05:19:01"
$ws.Range("N13").Value = "3.9 °C 4:58 TU"
$ws.Range("E14").Value = "2026-03-01 05:19:03"
$ws.Range("L14").Value = "10.1 km/h - 318º 4:34 TU"
$ws.Range("E15").Value = "2026-03-01 05:19:06"
$ws.Range("H15").Value = "87%"
$ws.Range("O15").Value = "8.4 °C"
$ws.Range("E16").Value = "2026-03-01 05:19:09"
$ws.Range("N16").Value = "-6.2 °C 4:36 TU"
$ws.Range("O16").Value = "-4.9 °C"
$ws.Range("E17").Value = "2026-03-01 05:19:11"
$ws.Range("L17").Value = "14.4 km/h - 224º 4:55 TU"
$ws.Range("E18").Value = "2026-03-01 05:19:14"
$ws.Range("N18").Value = "5.9 °C 4:59 TU"
$ws.Range("O18").Value = "7.0 °C"
$ws.Range("E19").Value = "2026-03-01 05:19:16"
$ws.Range("N19").Value = "5.9 °C 4:57 TU"
$ws.Range("E20").Value = "2026-03-01 05:19:19"
$ws.Range("O20").Value = "-3.0 °C"
$ws.Range("E21").Value = "2026-03-01 05:19:22"
$ws.Range("N21").Value = "6.0 °C 4:59 TU"
$ws.Range("E22").Value = "2026-03-01 05:19:24"
$ws.Range("L22").Value = "13.7 km/h - 331º 4:42 TU"
$ws.Range("E23").Value = "2026-03-01 05:19:27"
$ws.Range("N23").Value = "-4.2 °C 4:54 TU"
$ws.Range("O23").Value = "-3.5 °C"
$ws.Range("E24").Value = "2026-03-01 05:19:30"
$ws.Range("O24").Value = "4.4 °C"
$ws.Range("E25").Value = "2026-03-01 05:19:32"
$ws.Range("N25").Value = "-2.8 °C 4:30 TU"
$ws.Range("O25").Value = "-2.1 °C"
$ws.Range("E26").Value = "2026-03-01 05:19:35"
$ws.Range("H26").Value = "98%"
$ws.Range("N26").Value = "2.4 °C 4:34 TU"
$ws.Range("E27").Value = "2026-03-01 05:19:38"
$ws.Range("L27").Value = "9.7 km/h - 207º 4:53 TU"
$ws.Range("N27").Value = "-2.1 °C 4:59 TU"
$ws.Range("O27").Value = "-1.4 °C"
$ws.Range("E28").Value = "2026-03-01 05:19:41"
$ws.Range("N28").Value = "8.4 °C 4:59 TU"
$ws.Range("E29").Value = "2026-03-01 05:19:43"
$ws.Range("E30").Value = "2026-03-01 05:19:46"
$ws.Range("H30").Value = "78%"
$ws.Range("O30").Value = "10.2 °C"
$ws.Range("E31").Value = "2026-03-01 05:19:48"
$ws.Range("E32").Value = "2026-03-01 05:19:51"
$ws.Range("M32").Value = "4.5 °C 4:59 TU"
$ws.Range("O32").Value = "2.3 °C"
$ws.Range("E33").Value = "2026-03-01 05:19:54"
$ws.Range("E34").Value = "2026-03-01 05:19:57"
$ws.Range("L34").Value = "9.0 km/h - 169º 4:34 TU"
$ws.Range("N34").Value = "-0.5 °C 4:38 TU"
$ws.Range("E35").Value = "2026-03-01 05:19:59"
$ws.Range("E36").Value = "2026-03-01 05:20:02"
$ws.Range("L36").Value = "13.3 km/h - 18º 4:49 TU"
$ws.Range("M36").Value = "12.3 °C 4:59 TU"
$ws.Range("O36").Value = "9.7 °C"
$ws.Range("E37").Value = "2026-03-01 05:20:04"
$ws.Range("N37").Value = "6.1 °C 4:54 TU"
$ws.Range("E38").Value = "2026-03-01 05:20:07"
$ws.Range("E39").Value = "2026-03-01 05:20:09"
$ws.Range("E40").Value = "2026-03-01 05:20:12"
$ws.Range("N40").Value = "6.1 °C 4:59 TU"
$ws.Range("E41").Value = "2026-03-01 05:20:15"
$ws.Range("N41").Value = "11.2 °C 4:46 TU"
$ws.Range("E42").Value = "2026-03-01 05:20:17"
$ws.Range("H42").Value = "86%"
$ws.Range("N42").Value = "6.3 °C 4:33 TU"
$ws.Range("O42").Value = "8.7 °C"
$ws.Range("E43").Value = "2026-03-01 05:20:19"
$ws.Range("L43").Value = "5.8 km/h - 188º 4:34 TU"
$ws.Range("N43").Value = "8.2 °C 4:53 TU"
$ws.Range("E44").Value = "2026-03-01 05:20:22"
$ws.Range("N44").Value = "-3.3 °C 4:53 TU"
$ws.Range("O44").Value = "-2.6 °C"
$ws.Range("E45").Value = "2026-03-01 05:20:25"
$ws.Range("N45").Value = "2.8 °C 4:59 TU"
$ws.Range("E46").Value = "2026-03-01 05:20:27"
$ws.Range("J46").Value = "1026.5 hPa"
$ws.Range("O46").Value = "7.7 °C"
